$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mark the 4 newly developed routes as DONE (column I) for rows 6-9:
#   row 6 -> login
#   row 7 -> viewing a user's collections
#   row 8 -> viewing a collection's detail
#   row 9 -> adding a new collection
# Copy the formatting already used for the "DONE" mark in column I (I5, big
# centered font) onto the new cells, then set their value.
$ws.Range("I5").Copy()
$ws.Range("I6:I9").PasteSpecial(-4122)
$ws.Range("I6:I9").Value = "X"

# Update the active selection to the last cell worked on
$ws.Range("J10").Select()
